# Regenerate merged AHB files
#
# The sheet compares a "FV2310" (old) and "FV2404" (new) formatversion of
# an AHB (Anwendungshandbuch) segment table side by side. The header row
# used to carry generic "_old"/"_new" suffixes; relabel them with the
# concrete format-version tags, then present the data as a proper Excel
# Table with the header row frozen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) -------------------------------------
# Columns A:J described the "old" (FV2310) side, K is the unchanged "diff"
# marker column, and L:U described the "new" (FV2404) side.
$headerRenames = @{
  "A1" = "Segmentname_FV2310"
  "B1" = "Segmentgruppe_FV2310"
  "C1" = "Segment_FV2310"
  "D1" = "Datenelement_FV2310"
  "E1" = "Segment ID_FV2310"
  "F1" = "Code_FV2310"
  "G1" = "Qualifier_FV2310"
  "H1" = "Beschreibung_FV2310"
  "I1" = "Bedingungsausdruck_FV2310"
  "J1" = "Bedingung_FV2310"
  "L1" = "Segmentname_FV2404"
  "M1" = "Segmentgruppe_FV2404"
  "N1" = "Segment_FV2404"
  "O1" = "Datenelement_FV2404"
  "P1" = "Segment ID_FV2404"
  "Q1" = "Code_FV2404"
  "R1" = "Qualifier_FV2404"
  "S1" = "Beschreibung_FV2404"
  "T1" = "Bedingungsausdruck_FV2404"
  "U1" = "Bedingung_FV2404"
}

foreach ($addr in $headerRenames.Keys) {
  $ws.Range($addr).Value = $headerRenames[$addr]
}

# --- 2. Turn the data range into an Excel Table (ListObject) --------------
$tableRange = $ws.Range("A1:U68")
$listObject = $ws.ListObjects.Add(
  [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
  $tableRange,
  $null,
  [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# --- 3. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit applied"
